$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Forward/Reverse reads count + MD5 checksums for rows 3-10 (columns AB:AE)
$data = @(
    @(3,  22022227, 22022227, "359a090d1afe061bfc01b50dfa6e858b", "51ab0592e54e719c7331a46ab335399b"),
    @(4,  12526102, 12526102, "44e8e605b9d4a42a118a24441f87364a", "644415116fbd7953acbd098482fde289"),
    @(5,  25312891, 25312891, "d851d8dbfe370ac068659fddb9c21812", "d803a56b1766fa875dd037b760e8cad0"),
    @(6,  15542524, 15542524, "445760d27be328717da2515d23875f4b", "27689afea850d7376c89b82f90dac833"),
    @(7,  42326105, 42326105, "84382fb760336ac7d6315fbbf32b8c6f", "c757a8a59bc9cf04c51f0ab69118cbc7"),
    @(8,  12722007, 12722007, "2ee2460a6af40b41640b4f74d68a8268", "b6bb3c7cb72fa0bd52c20f1a99bd2976"),
    @(9,  20061376, 20061376, "f8fb2f8f18445b0fe7ebe1561e4795b1", "d5465756ea2aaf19339dad3b6b1f015a"),
    @(10, 17396379, 17396379, "834536a232530dcdcc29372d376deea9", "d60b7e31515bc7d8f54532c957a5497e")
)

# Fill numeric read counts first
foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("AB$r").Value = $row[1]
    $ws.Range("AC$r").Value = $row[2]
}

# Then Forward MD5 checksum column (AD) for all rows
foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("AD$r").Value = $row[3]
}

# Then Reverse MD5 checksum column (AE) for all rows
foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("AE$r").Value = $row[4]
}

# Scroll the frozen bottom-right pane so column X is the first visible
# column, and land the selection on AC13 (matches the saved view state)
$win = $excel.ActiveWindow
$win.ScrollColumn = 24
$ws.Range("AC13").Select()
